# Error handling for manager when wrong product details are entered
#
# Summary of the change applied to the "add products" template:
#  1. The "customizable" column (column E) is removed entirely - the
#     columns to its right (occasion, product_image, text_mask,
#     modal_mask) shift one place to the left.
#  2. Row 4 (TestProduct3)'s description is corrected from "testing3"
#     to "tesing 3".
#  3. A new product row is appended (TestProduct4 / testing4) with the
#     same price/category/occasion/image values as the other rows.
#  4. The worksheet view is reset (no frozen/scrolled topLeftCell, and
#     selection moved to B21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$imagePath = "/Users/leharbhatt/Desktop/samsung-galaxy-s10-black-back.png"

# --- 1. Remove the "customizable" column (column E) ---------------------
$ws.Range("E1").EntireColumn.Delete()

# --- 2. Fix the typo in TestProduct3's description (now row 4, col B) ---
$ws.Cells.Item(4, 2).Value = "tesing 3"

# --- 3. Append the new TestProduct4 row (row 5) --------------------------
$ws.Cells.Item(5, 1).Value = "TestProduct4"
$ws.Cells.Item(5, 2).Value = "testing4"
$ws.Cells.Item(5, 3).Value = 29.79
$ws.Cells.Item(5, 4).Value = "Placard"
$ws.Cells.Item(5, 5).Value = "Christmas"
$ws.Cells.Item(5, 6).Value = $imagePath
$ws.Cells.Item(5, 7).Value = $imagePath
$ws.Cells.Item(5, 8).Value = $imagePath

# Match the styling used by the other rows for the image columns (F:H)
$ws.Range("F2:H2").Copy()
$ws.Range("F5:H5").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- 4. Reset the worksheet view -----------------------------------------
$ws.Range("B21").Select()
